$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.440.96"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.980.92"
$ws.Range("E3").Value = "  +4.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9951"
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.00"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9963"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4636"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3940"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.11"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07907"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.45"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "1.980.23"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.193"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.843"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07094"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.71"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9986"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009937"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9968"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "29.508.87"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.555"
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "2.217.17"
$ws.Range("E25").Value = "  +9.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.111"
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.21"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.60"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.003"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.20"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.915"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09397"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.000004307"
$ws.Range("E33").Value = "  +162.76%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.8923"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.268"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.344"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.162"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05812"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.176"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02127"
$ws.Range("E40").Value = "  +2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.930"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9955"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5760"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1816"
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.805"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.13"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5375"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.185"
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.629"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06965"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.866"
$ws.Range("E51").Value = "  +0.52%  "
